$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 6851
$ws1.Range("F15").Value = 1110
$ws1.Range("F16").Value = 16290
$ws1.Range("F17").Value = 1605
$ws1.Range("F22").Value = 11417
$ws1.Range("F25").Value = 4498

# Sheet "全部类型" (All Types) updates - mirrors the same events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 6851
$ws4.Range("F17").Value = 1110
$ws4.Range("F18").Value = 16290
$ws4.Range("F19").Value = 1605
$ws4.Range("F26").Value = 11417
$ws4.Range("F29").Value = 4498
